# Auto-generated: update crypto price (D) and 1h-volume (E) columns
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.295.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.677.25'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5260'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.39%  '

$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2688'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.19%  '

$ws.Range("E9").Value = '  +1.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '

$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.699.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.513'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5769'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008468'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.323.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.916'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("E19").Value = '  +0.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.182'
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.801'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.03%  '

$ws.Range("E26").Value = '  +6.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06429'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.365'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.319'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.581'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.581'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.49%  '

$ws.Range("E33").Value = '  +1.93%  '

$ws.Range("E34").Value = '  +0.93%  '

$ws.Range("E35").Value = '  +2.03%  '

$ws.Range("E36").Value = '  +1.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.739'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.281'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.117.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.00%  '

$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8712'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.49%  '

$ws.Range("E42").Value = '  +0.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.827.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.30%  '

$ws.Range("E46").Value = '  +1.27%  '

$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.142'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05263'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.050'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.80%  '
